$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '29.523.40'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").Value = "'" + '1.851.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("D4").Value = "'" + '0.9994'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'" + '242.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("E6").Value = '  -0.96%  '

$ws.Range("D7").Value = "'" + '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = "'" + '47.87'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.49%  '

$ws.Range("D9").Value = "'" + '0.07570'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.12%  '

$ws.Range("D10").Value = "'" + '0.2978'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.47%  '

$ws.Range("D11").Value = "'" + '24.32'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").Value = "'" + '0.07683'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.41%  '

$ws.Range("D13").Value = "'" + '1.880.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.83%  '

$ws.Range("D14").Value = "'" + '5.029'
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'" + '0.6862'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").Value = "'" + '83.88'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.20%  '

$ws.Range("D17").Value = "'" + '0.000009878'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.34%  '

$ws.Range("D18").Value = "'" + '2.131.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.82%  '

$ws.Range("D19").Value = "'" + '6.213'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.80%  '

$ws.Range("D20").Value = "'" + '29.566.83'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.55%  '

$ws.Range("D21").Value = "'" + '235.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").Value = "'" + '12.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.89%  '

$ws.Range("D23").Value = "'" + '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("D24").Value = "'" + '7.633'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.38%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").Value = "'" + '155.82'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'" + '0.1387'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.81%  '

$ws.Range("D28").Value = "'" + '8.450'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.74%  '

$ws.Range("D29").Value = "'" + '17.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.86%  '

$ws.Range("D30").Value = "'" + '1.486'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.62%  '

$ws.Range("D31").Value = "'" + '0.05843'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.15%  '

$ws.Range("D32").Value = "'" + '1.278'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.33%  '

$ws.Range("D33").Value = "'" + '4.117'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.71%  '

$ws.Range("D34").Value = "'" + '4.052'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.84%  '

$ws.Range("E35").Value = '  +0.38%  '

$ws.Range("D36").Value = "'" + '1.171'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("D37").Value = "'" + '0.7178'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.28%  '

$ws.Range("D38").Value = "'" + '2.592'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.70%  '

$ws.Range("D39").Value = "'" + '2.797'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.18%  '

$ws.Range("D40").Value = "'" + '1.236.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.94%  '

$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("D42").Value = "'" + '0.9141'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.83%  '

$ws.Range("D43").Value = "'" + '6.147'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.23%  '

$ws.Range("D44").Value = "'" + '2.040.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.93%  '

$ws.Range("D45").Value = "'" + '0.9999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.06%  '

$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = "'" + '101.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.25%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = "'" + '67.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.92%  '

$ws.Range("D48").Value = "'" + '7.339'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.06%  '

$ws.Range("D49").Value = "'" + '9.177'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = "'" + '0.00000000117'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.16%  '

$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = "'" + '0.4040'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.27%  '

